$d = $word.ActiveDocument
$tbl = $d.Tables(1)

$tbl.Rows(1).Cells(1).Range.Text = "0M"
$tbl.Rows(2).Cells(1).Range.Text = "0M"
$tbl.Rows(3).Cells(1).Range.Text = "0M"
$tbl.Rows(4).Cells(1).Range.Text = "218"
$tbl.Rows(5).Cells(1).Range.Text = "0.00004"
$tbl.Rows(6).Cells(1).Range.Text = "0.00027"
$tbl.Rows(8).Cells(1).Range.Text = "0.00003"
$tbl.Rows(9).Cells(1).Range.Text = "0.00007"
$tbl.Rows(10).Cells(1).Range.Text = "0.00008"
$tbl.Rows(11).Cells(1).Range.Text = "0.00014"
$tbl.Rows(12).Cells(1).Range.Text = "0.01416"

$tbl.Rows(44).Cells(1).Range.Text = "99.99"
$tbl.Rows(45).Cells(1).Range.Text = "0.01"
$tbl.Rows(46).Cells(1).Range.Text = "265"
